$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 34
$ws.Range("B2").Value = 106
$ws.Range("B3").Value = 165
$ws.Range("B4").Value = 236
$ws.Range("B5").Value = 282
$ws.Range("B6").Value = 347
$ws.Range("B7").Value = 448
$ws.Range("B8").Value = 626
$ws.Range("B9").Value = 653
$ws.Range("B10").Value = 792
$ws.Range("B11").Value = 800
$ws.Range("B12").Value = 855
